# Updates cryptos list values per commit "Updated cryptos list on Thu Mar 21 19:26:05 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.131.26"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "'3.422.79"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'548.45"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").Value = "'177.77"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("D7").Value = "'0.642"
$ws.Range("E7").Value = "  +6.83%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  +6.17%  "
$ws.Range("D11").Value = "'53.32"
$ws.Range("E11").Value = "  -3.08%  "
$ws.Range("D12").Value = "'0.0000269"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").Value = "'9.13"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "'3.979.94"
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").Value = "'0.120"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "'3.421.01"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").Value = "'18.23"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "'65.150.70"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "'11.77"
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("D20").Value = "'0.979"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "'411.93"
$ws.Range("E21").Value = "  +5.82%  "
$ws.Range("D22").Value = "'3.99"
$ws.Range("E22").Value = "  +5.25%  "
$ws.Range("D23").Value = "'85.46"
$ws.Range("E23").Value = "  +3.22%  "
$ws.Range("D24").Value = "'4.10"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").Value = "'10.70"
$ws.Range("E25").Value = "  -10.11%  "
$ws.Range("D26").Value = "'2.83"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").Value = "'12.16"
$ws.Range("E27").Value = "  +5.32%  "
$ws.Range("D28").Value = "'6.02"
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("E29").Value = "  +5.03%  "
$ws.Range("D30").Value = "'29.65"
$ws.Range("E30").Value = "  +0.79%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'608.63"
$ws.Range("E31").Value = "  -8.16%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'6.46"
$ws.Range("E32").Value = "  -4.75%  "
$ws.Range("D33").Value = "'11.59"
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "'58.91"
$ws.Range("E35").Value = "  +2.08%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  +13.21%  "
$ws.Range("D38").Value = "'37.21"
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("D39").Value = "'0.0₃0774"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").Value = "'0.376"
$ws.Range("E40").Value = "  -5.47%  "
$ws.Range("D41").Value = "'3.179.91"
$ws.Range("E41").Value = "  +5.35%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "'2.52"
$ws.Range("E44").Value = "  -9.19%  "
$ws.Range("D45").Value = "'2.76"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").Value = "'3.24"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").Value = "'0.0408"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("E49").Value = "  +3.47%  "
$ws.Range("D50").Value = "'137.16"
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").Value = "'8.32"
$ws.Range("E51").Value = "  -2.06%  "
